$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Edit 1: row 553 (Grand Manan, 2005) -- TaxBaseCapita (G) recomputed ---
$ws.Range("G553").Value = 51.98573333333333

# --- Edit 2: row 554 previously duplicated Grand Manan's 2005 data; it should
# actually hold Grande-Anse's 2005 data (bugfix: Grande-Anse was dropped when
# parsing the 2005 tax base data and Grand Manan got duplicated in its place). ---
$ws.Range("B554").Value = "Grande-Anse"
$ws.Range("C554").Value = 1.341
$ws.Range("D554").Value = 0.09100000000000001
$ws.Range("E554").Value = 0.8130480656506448
$ws.Range("F554").Value = -0.688926143024619
$ws.Range("G554").Value = 36.42482415005862
$ws.Range("J554").Value = 853

# --- Edit 3: four municipalities (Bathurst, Campbellton, Dieppe, Edmundston)
# were incorrectly excluded when parsing the raw 2020 data. Insert four new
# rows at the top of the 2020 block (row 2018) and push the existing 2020
# rows down, then populate the new rows with the correct data. ---
$ws.Rows("2018:2021").Insert()

$ws.Range("A2018").Value = 2020
$ws.Range("B2018").Value = "Bathurst"
$ws.Range("C2018").Value = 1.775
$ws.Range("D2018").Value = 0.3607543918634951
$ws.Range("E2018").Value = 1.818663276456249
$ws.Range("F2018").Value = 0.5694793645456836
$ws.Range("G2018").Value = 90.7007474993696
$ws.Range("H2018").Value = $false
$ws.Range("I2018").Value = $false
$ws.Range("J2018").Value = 11897

$ws.Range("A2019").Value = 2020
$ws.Range("B2019").Value = "Campbellton"
$ws.Range("C2019").Value = 1.7763
$ws.Range("D2019").Value = 0.3826215313090222
$ws.Range("E2019").Value = 1.774466947551939
$ws.Range("F2019").Value = 0.6045397355804153
$ws.Range("G2019").Value = 87.40352317303503
$ws.Range("H2019").Value = $false
$ws.Range("I2019").Value = $true
$ws.Range("J2019").Value = 6883

$ws.Range("A2020").Value = 2020
$ws.Range("B2020").Value = "Dieppe"
$ws.Range("C2020").Value = 1.6295
$ws.Range("D2020").Value = 0.2344724235739048
$ws.Range("E2020").Value = 2.087807004412228
$ws.Range("F2020").Value = 0.2216123148439962
$ws.Range("G2020").Value = 128.9160180822566
$ws.Range("H2020").Value = $false
$ws.Range("I2020").Value = $true
$ws.Range("J2020").Value = 25384

$ws.Range("A2021").Value = 2020
$ws.Range("B2021").Value = "Edmundston"
$ws.Range("C2021").Value = 1.635
$ws.Range("D2021").Value = 0.33544794933655
$ws.Range("E2021").Value = 1.806330398069964
$ws.Range("F2021").Value = 0.7111211097708079
$ws.Range("G2021").Value = 87.5019722557298
$ws.Range("H2021").Value = $false
$ws.Range("I2021").Value = $false
$ws.Range("J2021").Value = 16580

# --- Edit 4: grow the "Frame0" table to cover the 4 newly-inserted rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J2110"))
